$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Columns("F:F").Insert()
$ws.Cells.Item(1, 6).Value = "Plate_Barcode"
$ws.Range("F7").Select()
